$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Check function..." / "GetValidationReportDS" row (row 3) entirely
# without shifting the rows below it: delete the two cells, then delete the
# (now-empty) row and re-insert a blank row in its place so every other row
# keeps its original row number.
$ws.Range("B3:C3").Delete()
$ws.Rows(3).Insert()
$ws.Rows(3).ClearFormats()
$ws.Rows(3).ClearContents()

# Update the selection to match the new focal cell
[void]$ws.Range("B6").Select()
